$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 474.3
$ws.Cells.Item(2, 9).Value = 517.875
$ws.Cells.Item(2, 11).Value = 517.875
$ws.Cells.Item(2, 13).Value = -404.875

$ws.Cells.Item(38, 8).Value = 365.7143
$ws.Cells.Item(38, 9).Value = 93.333336
$ws.Cells.Item(38, 10).Value = 2000
$ws.Cells.Item(38, 11).Value = 280.000008
$ws.Cells.Item(38, 12).Value = 6000
$ws.Cells.Item(38, 13).Value = 91.99999200000002
$ws.Cells.Item(38, 14).Value = -6744

$ws.Cells.Item(43, 8).Value = 1549.6666
$ws.Cells.Item(43, 9).Value = 1191.1666
$ws.Cells.Item(43, 11).Value = 1191.1666
$ws.Cells.Item(43, 13).Value = -1122.1666

$ws.Cells.Item(58, 8).Value = 2828.125
$ws.Cells.Item(58, 9).Value = 325
$ws.Cells.Item(58, 10).Value = 7000
$ws.Cells.Item(58, 11).Value = 975
$ws.Cells.Item(58, 12).Value = 21000
$ws.Cells.Item(58, 13).Value = -825
$ws.Cells.Item(58, 14).Value = -21300

$ws.Cells.Item(132, 8).Value = 1384.2264
$ws.Cells.Item(132, 9).Value = 1122.2693
$ws.Cells.Item(132, 11).Value = 3366.8079
$ws.Cells.Item(132, 13).Value = -836.8078999999998

$ws.Cells.Item(137, 8).Value = 560226
$ws.Cells.Item(137, 9).Value = 1507.625
$ws.Cells.Item(137, 10).Value = 1454175.4
$ws.Cells.Item(137, 11).Value = 4522.875
$ws.Cells.Item(137, 12).Value = 4362526.199999999
$ws.Cells.Item(137, 13).Value = -1972.875
$ws.Cells.Item(137, 14).Value = -4367626.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4477.507
$ws.Cells.Item(32, 9).Value = 1891.569
$ws.Cells.Item(32, 10).Value = 14476.467
$ws.Cells.Item(32, 11).Value = 1891.569
$ws.Cells.Item(32, 12).Value = 14476.467
$ws.Cells.Item(32, 13).Value = -1604.569
$ws.Cells.Item(32, 14).Value = -15050.467

$ws.Cells.Item(50, 8).Value = 1134.2222
$ws.Cells.Item(50, 9).Value = 1524.1666
$ws.Cells.Item(50, 11).Value = 1524.1666
$ws.Cells.Item(50, 13).Value = -810.1666

$ws.Cells.Item(132, 8).Value = 3077.8545
$ws.Cells.Item(132, 9).Value = 2507.8262
$ws.Cells.Item(132, 11).Value = 7523.4786
$ws.Cells.Item(132, 13).Value = -4993.4786

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 62147.766
$ws.Cells.Item(105, 9).Value = 74815.14
$ws.Cells.Item(105, 10).Value = 3033.3333
$ws.Cells.Item(105, 11).Value = 74815.14
$ws.Cells.Item(105, 12).Value = 3033.3333
$ws.Cells.Item(105, 13).Value = -73068.14
$ws.Cells.Item(105, 14).Value = -6527.3333

$ws.Cells.Item(134, 8).Value = 3363.5112
$ws.Cells.Item(134, 9).Value = 2260.5293
$ws.Cells.Item(134, 10).Value = 6772.727
$ws.Cells.Item(134, 11).Value = 6781.5879
$ws.Cells.Item(134, 12).Value = 20318.181
$ws.Cells.Item(134, 13).Value = -4246.5879
$ws.Cells.Item(134, 14).Value = -25388.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1802.2325
$ws.Cells.Item(31, 10).Value = 2175.52
$ws.Cells.Item(31, 12).Value = 2175.52
$ws.Cells.Item(31, 14).Value = -2765.52

$ws.Cells.Item(34, 8).Value = 1802.2325
$ws.Cells.Item(34, 10).Value = 2175.52
$ws.Cells.Item(34, 12).Value = 2175.52
$ws.Cells.Item(34, 14).Value = -2579.52

$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).Value = $null

$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).Value = $null

$ws.Cells.Item(102, 8).Value = 40600
$ws.Cells.Item(102, 10).Value = 40600
$ws.Cells.Item(102, 12).Value = 40600
$ws.Cells.Item(102, 14).Value = -45468

$ws.Cells.Item(122, 8).Value = 2262.1667
$ws.Cells.Item(122, 9).Value = 1508.0769
$ws.Cells.Item(122, 10).Value = 2838.8235
$ws.Cells.Item(122, 11).Value = 4524.2307
$ws.Cells.Item(122, 12).Value = 8516.470499999999
$ws.Cells.Item(122, 13).Value = -2074.2307
$ws.Cells.Item(122, 14).Value = -13416.4705

$ws.Cells.Item(132, 8).Value = 1952497
$ws.Cells.Item(132, 9).Value = 1980804.4
$ws.Cells.Item(132, 11).Value = 5942413.199999999
$ws.Cells.Item(132, 13).Value = -5939883.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1118.8823
$ws.Cells.Item(107, 10).Value = 1167.8572
$ws.Cells.Item(107, 12).Value = 3503.5716
$ws.Cells.Item(107, 14).Value = -7343.571599999999

$ws.Cells.Item(122, 8).Value = 415.77777
$ws.Cells.Item(122, 10).Value = 422.57144
$ws.Cells.Item(122, 12).Value = 3803.14296
$ws.Cells.Item(122, 14).Value = -8703.142960000001

$ws.Cells.Item(131, 8).Value = 1091.9048
$ws.Cells.Item(131, 10).Value = 1357.8182
$ws.Cells.Item(131, 12).Value = 4073.4546
$ws.Cells.Item(131, 14).Value = -14153.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 764.8461
$ws.Cells.Item(107, 9).Value = 813.0909
$ws.Cells.Item(107, 10).Value = 499.5
$ws.Cells.Item(107, 11).Value = 813.0909
$ws.Cells.Item(107, 12).Value = 499.5
$ws.Cells.Item(107, 13).Value = 1106.9091
$ws.Cells.Item(107, 14).Value = -4339.5

$ws.Cells.Item(113, 8).Value = 3868.5
$ws.Cells.Item(113, 9).Value = 3237
$ws.Cells.Item(113, 10).Value = 4815.75
$ws.Cells.Item(113, 11).Value = 3237
$ws.Cells.Item(113, 12).Value = 4815.75
$ws.Cells.Item(113, 13).Value = -1067
$ws.Cells.Item(113, 14).Value = -9155.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1152.8572
$ws.Cells.Item(16, 9).Value = 997.8889
$ws.Cells.Item(16, 10).Value = 2082.6667
$ws.Cells.Item(16, 11).Value = 997.8889
$ws.Cells.Item(16, 12).Value = 2082.6667
$ws.Cells.Item(16, 13).Value = -827.8889
$ws.Cells.Item(16, 14).Value = -2422.6667

$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).Value = $null

$ws.Cells.Item(61, 8).Value = 7336
$ws.Cells.Item(61, 9).Value = 7336
$ws.Cells.Item(61, 11).Value = 7336
$ws.Cells.Item(61, 13).Value = -7134

$ws.Cells.Item(63, 8).Value = 30000
$ws.Cells.Item(63, 9).Value = 20000
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 13).Value = -19251

$ws.Cells.Item(66, 8).Value = 30000
$ws.Cells.Item(66, 9).Value = 20000
$ws.Cells.Item(66, 11).Value = 60000
$ws.Cells.Item(66, 13).Value = -56256

$ws.Cells.Item(68, 8).Value = 2443.2856
$ws.Cells.Item(68, 9).Value = 2033.3334
$ws.Cells.Item(68, 11).Value = 2033.3334
$ws.Cells.Item(68, 13).Value = -1284.3334

$ws.Cells.Item(71, 8).Value = 2443.2856
$ws.Cells.Item(71, 9).Value = 2033.3334
$ws.Cells.Item(71, 11).Value = 10166.667
$ws.Cells.Item(71, 13).Value = -6422.666999999999

$ws.Cells.Item(93, 8).Value = 2638
$ws.Cells.Item(93, 9).Value = 2408.6365
$ws.Cells.Item(93, 11).Value = 2408.6365
$ws.Cells.Item(93, 13).Value = -1160.6365

$ws.Cells.Item(99, 8).Value = 44996.332
$ws.Cells.Item(99, 9).Value = 44996.332
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 44996.332
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -42001.332
$ws.Cells.Item(99, 14).Value = $null

$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).Value = $null

$ws.Cells.Item(113, 8).Value = 7336
$ws.Cells.Item(113, 9).Value = 7336
$ws.Cells.Item(113, 11).Value = 7336
$ws.Cells.Item(113, 13).Value = -5166

$ws.Cells.Item(117, 8).Value = 50502.223
$ws.Cells.Item(117, 10).Value = 50502.223
$ws.Cells.Item(117, 12).Value = 50502.223
$ws.Cells.Item(117, 14).Value = -59680.223

$ws.Cells.Item(122, 8).Value = 71575360
$ws.Cells.Item(122, 9).Value = 91094720
$ws.Cells.Item(122, 10).Value = 4368
$ws.Cells.Item(122, 11).Value = 273284160
$ws.Cells.Item(122, 12).Value = 13104
$ws.Cells.Item(122, 13).Value = -273281710
$ws.Cells.Item(122, 14).Value = -18004

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 61445
$ws.Cells.Item(27, 10).Value = 61445
$ws.Cells.Item(27, 12).Value = 61445
$ws.Cells.Item(27, 14).Value = -61583

$ws.Cells.Item(33, 8).Value = 23000
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 23000
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 23000
$ws.Cells.Item(33, 13).Value = $null
$ws.Cells.Item(33, 14).Value = -23500

$ws.Cells.Item(36, 8).Value = 23000
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 23000
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 23000
$ws.Cells.Item(36, 13).Value = $null
$ws.Cells.Item(36, 14).Value = -23500

$ws.Cells.Item(102, 8).Value = 47000
$ws.Cells.Item(102, 10).Value = 47000
$ws.Cells.Item(102, 12).Value = 47000
$ws.Cells.Item(102, 14).Value = -53490

$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).Value = $null

$ws.Cells.Item(122, 8).Value = 2470.6333
$ws.Cells.Item(122, 9).Value = 2440.76
$ws.Cells.Item(122, 10).Value = 2620
$ws.Cells.Item(122, 11).Value = 7322.280000000001
$ws.Cells.Item(122, 12).Value = 7860
$ws.Cells.Item(122, 13).Value = -4872.280000000001
$ws.Cells.Item(122, 14).Value = -12760

$ws.Cells.Item(126, 8).Value = 3583.238
$ws.Cells.Item(126, 9).Value = 3006.6365
$ws.Cells.Item(126, 11).Value = 9019.9095
$ws.Cells.Item(126, 13).Value = -6549.9095

$ws.Cells.Item(132, 8).Value = 2572.3809
$ws.Cells.Item(132, 9).Value = 2468.1516
$ws.Cells.Item(132, 10).Value = 2954.5557
$ws.Cells.Item(132, 11).Value = 7404.4548
$ws.Cells.Item(132, 12).Value = 8863.667099999999
$ws.Cells.Item(132, 13).Value = -4874.4548
$ws.Cells.Item(132, 14).Value = -13923.6671

$ws.Cells.Item(136, 8).Value = 1817.6072
$ws.Cells.Item(136, 9).Value = 1527.9565
$ws.Cells.Item(136, 10).Value = 3150
$ws.Cells.Item(136, 11).Value = 4583.8695
$ws.Cells.Item(136, 12).Value = 9450
$ws.Cells.Item(136, 13).Value = -2033.8695
$ws.Cells.Item(136, 14).Value = -14550
